$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 128.57143
$ws.Range("I9").Value = 120
$ws.Range("J9").Value = 150
$ws.Range("K9").Value = 120
$ws.Range("L9").Value = 150
$ws.Range("M9").Value = 49
$ws.Range("N9").Value = -488

# Row 12
$ws.Range("H12").Value = 97.333336
$ws.Range("I12").Value = 120
$ws.Range("J12").Value = 86
$ws.Range("K12").Value = 120
$ws.Range("L12").Value = 86
$ws.Range("M12").Value = 50
$ws.Range("N12").Value = -426

# Row 18
$ws.Range("H18").Value = 400
$ws.Range("I18").Value = 400
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 400
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -116

# Row 21
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("N21").ClearContents()

# Row 23
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("N23").ClearContents()

# Row 29
$ws.Range("H29").Value = 812.5
$ws.Range("I29").Value = 700
$ws.Range("J29").Value = 1000
$ws.Range("K29").Value = 2100
$ws.Range("L29").Value = 3000
$ws.Range("M29").Value = -1819
$ws.Range("N29").Value = -3562

# Row 38
$ws.Range("H38").Value = 347.5
$ws.Range("I38").Value = 86.38461
$ws.Range("J38").Value = 832.4286
$ws.Range("K38").Value = 259.15383
$ws.Range("L38").Value = 2497.2858
$ws.Range("M38").Value = 112.84617
$ws.Range("N38").Value = -3241.2858

# Row 40
$ws.Range("H40").Value = 1797.3334
$ws.Range("I40").Value = 1200
$ws.Range("J40").Value = 2014.5454
$ws.Range("K40").Value = 1200
$ws.Range("L40").Value = 2014.5454
$ws.Range("M40").Value = -1025
$ws.Range("N40").Value = -2364.5454

# Row 58
$ws.Range("H58").Value = 951.88464
$ws.Range("I58").Value = 435.7143
$ws.Range("J58").Value = 1142.0526
$ws.Range("K58").Value = 1307.1429
$ws.Range("L58").Value = 3426.1578
$ws.Range("M58").Value = -1157.1429
$ws.Range("N58").Value = -3726.1578

# Row 87
$ws.Range("H87").Value = 35095.715
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 35095.715
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 35095.715
$ws.Range("N87").Value = -37591.715

# Row 90
$ws.Range("H90").Value = 35095.715
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 35095.715
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 105287.145
$ws.Range("N90").Value = -117767.145

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 511964.47
$ws.Range("I32").Value = 3245.5103
$ws.Range("J32").Value = 3004687.5
$ws.Range("K32").Value = 3245.5103
$ws.Range("L32").Value = 3004687.5
$ws.Range("M32").Value = -2958.5103
$ws.Range("N32").Value = -3005261.5

$ws = $wb.Worksheets.Item("CRP")
# Row 17
$ws.Range("H17").Value = 100000000
$ws.Range("I17").Value = 100000000
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 100000000
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -99999826

# Row 31
$ws.Range("H31").Value = 13415.756
$ws.Range("I31").Value = 4221.067
$ws.Range("J31").Value = 38492.184
$ws.Range("K31").Value = 4221.067
$ws.Range("L31").Value = 38492.184
$ws.Range("M31").Value = -3926.067
$ws.Range("N31").Value = -39082.184

# Row 34
$ws.Range("H34").Value = 13415.756
$ws.Range("I34").Value = 4221.067
$ws.Range("J34").Value = 38492.184
$ws.Range("K34").Value = 4221.067
$ws.Range("L34").Value = 38492.184
$ws.Range("M34").Value = -4019.067
$ws.Range("N34").Value = -38896.184

# Row 41
$ws.Range("H41").Value = 100000000
$ws.Range("I41").Value = 100000000
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 100000000
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -99999572
$ws.Range("N41").ClearContents()

# Row 50
$ws.Range("H50").Value = 11635.857
$ws.Range("I50").Value = 9083
$ws.Range("J50").Value = 12061.333
$ws.Range("K50").Value = 9083
$ws.Range("L50").Value = 12061.333
$ws.Range("M50").Value = -8458
$ws.Range("N50").Value = -13311.333

# Row 51
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("N51").ClearContents()

# Row 59
$ws.Range("H59").Value = 31460.334
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 31460.334
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 31460.334
$ws.Range("N59").Value = -33750.334

# Row 60
$ws.Range("H60").Value = 12306.823
$ws.Range("I60").Value = 2000
$ws.Range("J60").Value = 13681.066
$ws.Range("K60").Value = 2000
$ws.Range("L60").Value = 13681.066
$ws.Range("M60").Value = -1489
$ws.Range("N60").Value = -14703.066

# Row 61
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").ClearContents()

# Row 62
$ws.Range("I62").Value = 2800
$ws.Range("J62").Value = 3333.3333
$ws.Range("K62").Value = 2800
$ws.Range("L62").Value = 3333.3333
$ws.Range("M62").Value = -2176
$ws.Range("N62").Value = -4581.3333

# Row 65
$ws.Range("I65").Value = 2800
$ws.Range("J65").Value = 3333.3333
$ws.Range("K65").Value = 14000
$ws.Range("L65").Value = 16666.6665
$ws.Range("M65").Value = -10880
$ws.Range("N65").Value = -22906.6665

# Row 68
$ws.Range("H68").Value = 25483.75
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 25483.75
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 25483.75
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -26981.75

# Row 71
$ws.Range("H71").Value = 25483.75
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 25483.75
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 76451.25
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -83939.25

# Row 74
$ws.Range("H74").Value = 15898.875
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 15898.875
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 15898.875
$ws.Range("N74").Value = -17646.875

# Row 77
$ws.Range("H77").Value = 15898.875
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 15898.875
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 47696.625
$ws.Range("N77").Value = -56432.625

$ws = $wb.Worksheets.Item("CUL")
# Row 132
$ws.Range("H132").Value = 1067.8846
$ws.Range("I132").Value = 666.6667
$ws.Range("J132").Value = 1501.2
$ws.Range("K132").Value = 6000.0003
$ws.Range("L132").Value = 13510.8
$ws.Range("M132").Value = -3470.0003
$ws.Range("N132").Value = -18570.8

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 20836734
$ws.Range("I81").Value = 1971.5714
$ws.Range("J81").Value = 50005400
$ws.Range("K81").Value = 3943.1428
$ws.Range("L81").Value = 100010800
$ws.Range("M81").Value = -2882.1428
$ws.Range("N81").Value = -100012922

# Row 84
$ws.Range("H84").Value = 20836734
$ws.Range("I84").Value = 1971.5714
$ws.Range("J84").Value = 50005400
$ws.Range("K84").Value = 19715.714
$ws.Range("L84").Value = 500054000
$ws.Range("M84").Value = -14411.714
$ws.Range("N84").Value = -500064608
